$wb = $excel.ActiveWorkbook

# --- Work on the "Repayment Schedule" sheet (sheet4.xml) ---
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before column N, shifting existing N:P (Late/blank/Outstanding)
# to O:Q. This matches the dimension growing from A1:P15 to A1:Q15.
$wsSchedule.Columns("N").Insert()

# Make "Repayment Schedule" the active tab (previously "Transactions" was active),
# and move its selection to S6.
$wsSchedule.Activate()
$wsSchedule.Range("S6").Select()
